# "Add more task on Sprint Backlog.xlsx Update Group 2 Wiki.docx"
#
# Project Backlog sheet: bump two estimate values, add Sprint(=2) tags,
# add a new backlog row (#15), widen the table to include it, and
# tighten column C's width (Excel "best fit").
#
# Impediment Backlog sheet: log a new impediment (row 5) with the same
# "Not solve" status styling as the rows above it, and widen column C
# to fit the new (longer) text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Project Backlog")
$ws2 = $wb.Worksheets.Item("Impediment Backlog")

# --- Project Backlog ---------------------------------------------------

# Re-estimate rows 4 & 5 (value 10 -> 200) and tag them Sprint 2.
$ws1.Range("E4").Value = 200
$ws1.Range("G4").Value = 2
$ws1.Range("E5").Value = 200
$ws1.Range("G5").Value = 2

# New backlog item #15 on row 18.
$ws1.Range("A18").Value = 15

# Grow Table1 so the new row is part of it (and its AutoFilter).
$lo1 = $ws1.ListObjects.Item("Table1")
$lo1.Resize($ws1.Range("A2:G18"))

# Column C was manually widened before; now it's auto-fit back down.
$ws1.Columns.Item(3).ColumnWidth = 30.140625

# --- Impediment Backlog --------------------------------------------------

# New impediment row, mirroring the style of the existing "Not solve" rows.
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "Not solve"
$ws2.Range("C5").Value = "When finish task, not commit to SVN for other member verify it"
$ws2.Range("B4").Copy()
$ws2.Range("B5").PasteSpecial(-4122)

# Column C widened to fit the new (longer) impediment text.
$ws2.Columns.Item(3).ColumnWidth = 58.42578125

# --- View state: selections (sheet1 first, sheet2 last so Impediment
#     Backlog remains the active/visible tab, as in the original file) ---
$ws1.Range("B4").Select()
$ws2.Range("C6").Select()
